# Add a new attendance-date column (U) for 2025-11-22:
#  - U1 header, styled like the other date headers
#  - U2:U9 marked absent ("\u274c") for every student
#  - S2:S9 (Total) bumped by 1 to account for the new tracked day

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell U1 -------------------------------------------------
# Write as literal text (not an auto-parsed date) by temporarily forcing
# a text number format, then restore "General" by copying T1's format
# (which also carries over the bold/border/centered header style).
$ws.Range("U1").NumberFormat = "@"
$ws.Range("U1").Value = "2025-11-22"

$ws.Range("T1").Copy()
$ws.Range("U1").PasteSpecial(-4122)   # xlPasteFormats

# --- Data rows 2..9 ---------------------------------------------------
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 21).Value = "❌"   # column U = 21

    $total = $ws.Cells.Item($r, 19)       # column S = 19 ("Total")
    $current = $total.Value2
    $total.Value2 = $current + 1
}
